$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 6864629
$ws.Range("C36").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D36").Value = 45186.61458333334
$ws.Range("E36").Value = "Borac Banja Luka"
$ws.Range("F36").Value = "NK Posusje"
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = "H"
$ws.Range("L36").Value = 1.363
$ws.Range("M36").Value = 4.5
$ws.Range("N36").Value = 6.5
$ws.Range("O36").Value = 1.363
$ws.Range("P36").Value = 4.2
$ws.Range("Q36").Value = 6.5
$ws.Range("R36").Value = -1.25
$ws.Range("S36").Value = 1.95
$ws.Range("T36").Value = 1.85
$ws.Range("U36").Value = 2.5
$ws.Range("V36").Value = 1.925
$ws.Range("W36").Value = 1.875
$ws.Range("X36").Value = 0.363
$ws.Range("Y36").Value = -1
$ws.Range("Z36").Value = -1
$ws.Range("AA36").Value = -0.5
$ws.Range("AB36").Value = 0.425
$ws.Range("AC36").Value = -1
$ws.Range("AD36").Value = 0.875

# Row 37
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 6865299
$ws.Range("C37").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D37").Value = 45186.61458333334
$ws.Range("E37").Value = "Siroki Brijeg"
$ws.Range("F37").Value = "Zvijezda 09"
$ws.Range("G37").Value = 2
$ws.Range("H37").Value = 1
$ws.Range("I37").Value = 2
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = "H"
$ws.Range("L37").Value = 1.25
$ws.Range("M37").Value = 5.5
$ws.Range("N37").Value = 8
$ws.Range("O37").Value = 1.4
$ws.Range("P37").Value = 4.75
$ws.Range("Q37").Value = 5.75
$ws.Range("R37").Value = -1.25
$ws.Range("S37").Value = 1.9
$ws.Range("T37").Value = 1.9
$ws.Range("U37").Value = 2.75
$ws.Range("V37").Value = 1.85
$ws.Range("W37").Value = 1.95
$ws.Range("X37").Value = 0.3999999999999999
$ws.Range("Y37").Value = -1
$ws.Range("Z37").Value = -1
$ws.Range("AA37").Value = -0.5
$ws.Range("AB37").Value = 0.45
$ws.Range("AC37").Value = 0.425
$ws.Range("AD37").Value = -0.5

# Row 189
$ws.Range("A189").Value = 187
$ws.Range("B189").Value = 7952780
$ws.Range("C189").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D189").Value = 45432.5
$ws.Range("E189").Value = "Velez Mostar"
$ws.Range("F189").Value = "GOSK Gabela"
$ws.Range("G189").Value = 3
$ws.Range("H189").Value = 3
$ws.Range("I189").Value = 1
$ws.Range("J189").Value = 1
$ws.Range("K189").Value = "D"
$ws.Range("L189").Value = 1.4
$ws.Range("M189").Value = 4
$ws.Range("N189").Value = 7
$ws.Range("O189").Value = 1.363
$ws.Range("P189").Value = 4.2
$ws.Range("Q189").Value = 8
$ws.Range("R189").Value = -1.5
$ws.Range("S189").Value = 2
$ws.Range("T189").Value = 1.8
$ws.Range("U189").Value = 2.75
$ws.Range("V189").Value = 1.825
$ws.Range("W189").Value = 1.975
$ws.Range("X189").Value = -1
$ws.Range("Y189").Value = 3.2
$ws.Range("Z189").Value = -1
$ws.Range("AA189").Value = -1
$ws.Range("AB189").Value = 0.8
$ws.Range("AC189").Value = 0.825
$ws.Range("AD189").Value = -1

# Row 190
$ws.Range("A190").Value = 188
$ws.Range("B190").Value = 7952779
$ws.Range("C190").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D190").Value = 45432.5
$ws.Range("E190").Value = "Zrinjski Mostar"
$ws.Range("F190").Value = "FK Tuzla City"
$ws.Range("G190").Value = 4
$ws.Range("H190").Value = 0
$ws.Range("I190").Value = 2
$ws.Range("J190").Value = 0
$ws.Range("K190").Value = "H"
$ws.Range("L190").Value = 1.25
$ws.Range("M190").Value = 5.75
$ws.Range("N190").Value = 7
$ws.Range("O190").Value = 1.055
$ws.Range("P190").Value = 13
$ws.Range("Q190").Value = 17
$ws.Range("R190").Value = -3.5
$ws.Range("S190").Value = 1.975
$ws.Range("T190").Value = 1.825
$ws.Range("U190").Value = 4.75
$ws.Range("V190").Value = 1.825
$ws.Range("W190").Value = 1.975
$ws.Range("X190").Value = 0.05499999999999994
$ws.Range("Y190").Value = -1
$ws.Range("Z190").Value = -1
$ws.Range("AA190").Value = 0.9750000000000001
$ws.Range("AB190").Value = -1
$ws.Range("AC190").Value = -1
$ws.Range("AD190").Value = 0.9750000000000001

# Row 191
$ws.Range("A191").Value = 189
$ws.Range("B191").Value = 7952777
$ws.Range("C191").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D191").Value = 45432.5
$ws.Range("E191").Value = "Borac Banja Luka"
$ws.Range("F191").Value = "NK Igman Konjic"
$ws.Range("G191").Value = 4
$ws.Range("H191").Value = 3
$ws.Range("I191").Value = 1
$ws.Range("J191").Value = 2
$ws.Range("K191").Value = "H"
$ws.Range("L191").Value = 1.25
$ws.Range("M191").Value = 5.75
$ws.Range("N191").Value = 7
$ws.Range("O191").Value = 1.2
$ws.Range("P191").Value = 5.75
$ws.Range("Q191").Value = 12
$ws.Range("R191").Value = -2
$ws.Range("S191").Value = 1.95
$ws.Range("T191").Value = 1.85
$ws.Range("U191").Value = 3.25
$ws.Range("V191").Value = 1.9
$ws.Range("W191").Value = 1.9
$ws.Range("X191").Value = 0.2
$ws.Range("Y191").Value = -1
$ws.Range("Z191").Value = -1
$ws.Range("AA191").Value = -1
$ws.Range("AB191").Value = 0.8500000000000001
$ws.Range("AC191").Value = 0.8999999999999999
$ws.Range("AD191").Value = -1

# Row 197
$ws.Range("A197").Value = 195
$ws.Range("B197").Value = 7952787
$ws.Range("C197").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D197").Value = 45438.5
$ws.Range("E197").Value = "Zeljeznicar"
$ws.Range("F197").Value = "Borac Banja Luka"
$ws.Range("G197").Value = 2
$ws.Range("H197").Value = 1
$ws.Range("K197").Value = "H"
$ws.Range("L197").Value = 3.3
$ws.Range("M197").Value = 3
$ws.Range("N197").Value = 2.05
$ws.Range("O197").Value = 1.727
$ws.Range("P197").Value = 3.1
$ws.Range("Q197").Value = 4.5
$ws.Range("R197").Value = -0.5
$ws.Range("S197").Value = 1.825
$ws.Range("T197").Value = 1.975
$ws.Range("U197").Value = 2.25
$ws.Range("V197").Value = 2
$ws.Range("W197").Value = 1.8
$ws.Range("X197").Value = 0.7270000000000001
$ws.Range("Y197").Value = -1
$ws.Range("Z197").Value = -1
$ws.Range("AA197").Value = 0.825
$ws.Range("AB197").Value = -1
$ws.Range("AC197").Value = 1
$ws.Range("AD197").Value = -1

# Row 198
$ws.Range("A198").Value = 196
$ws.Range("B198").Value = 8259814
$ws.Range("C198").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D198").Value = 45438.5
$ws.Range("E198").Value = "Siroki Brijeg"
$ws.Range("F198").Value = "FK Sarajevo"
$ws.Range("G198").Value = 2
$ws.Range("H198").Value = 2
$ws.Range("K198").Value = "D"
$ws.Range("L198").Value = 3.4
$ws.Range("M198").Value = 3.1
$ws.Range("N198").Value = 2
$ws.Range("O198").Value = 9.5
$ws.Range("P198").Value = 4.75
$ws.Range("Q198").Value = 1.25
$ws.Range("R198").Value = 1.5
$ws.Range("S198").Value = 1.975
$ws.Range("T198").Value = 1.825
$ws.Range("U198").Value = 2.75
$ws.Range("V198").Value = 1.95
$ws.Range("W198").Value = 1.85
$ws.Range("X198").Value = -1
$ws.Range("Y198").Value = 3.75
$ws.Range("Z198").Value = -1
$ws.Range("AA198").Value = 0.9750000000000001
$ws.Range("AB198").Value = -1
$ws.Range("AC198").Value = 0.95
$ws.Range("AD198").Value = -1

# Row 199
$ws.Range("A199").Value = 197
$ws.Range("B199").Value = 8259815
$ws.Range("C199").Value = "Bosnia Herzegovina Premier Liga"
$ws.Range("D199").Value = 45438.5
$ws.Range("E199").Value = "NK Posusje"
$ws.Range("F199").Value = "Zvijezda 09"
$ws.Range("G199").Value = 2
$ws.Range("H199").Value = 0
$ws.Range("K199").Value = "H"
$ws.Range("L199").Value = 1.4
$ws.Range("M199").Value = 4
$ws.Range("N199").Value = 6.5
$ws.Range("O199").Value = 1.25
$ws.Range("P199").Value = 5
$ws.Range("Q199").Value = 8.5
$ws.Range("R199").Value = -1.75
$ws.Range("S199").Value = 2
$ws.Range("T199").Value = 1.8
$ws.Range("U199").Value = 3
$ws.Range("V199").Value = 1.925
$ws.Range("W199").Value = 1.875
$ws.Range("X199").Value = 0.25
$ws.Range("Y199").Value = -1
$ws.Range("Z199").Value = -1
$ws.Range("AA199").Value = 0.5
$ws.Range("AB199").Value = -0.5
$ws.Range("AC199").Value = -1
$ws.Range("AD199").Value = 0.875
